$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$grp = $s.Shapes.Item("Group 186")
$tb = $grp.GroupItems.Item("TextBox 188")

$tb.Height = 19.114882469763778
Write-Output "Height after: $($tb.Height)"
